$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values for Npnt-Itgb1 LR-pair analysis (YoungD7)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.824961333333333
$ws.Range("H2").Value = 5.474884
$ws.Range("I2").Value = 0.377840167393297
$ws.Range("J2").Value = 0.3778401673932969
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 306.7939346366765
$ws.Range("R2").Value = 2761.145411730088
$ws.Range("S2").Value = 0.1127543923907114
$ws.Range("T2").Value = 0.1127543923907114
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.824961333333333
$ws.Range("H3").Value = 5.474884
$ws.Range("I3").Value = 0.377840167393297
$ws.Range("J3").Value = 0.3778401673932969
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 297.4800771838876
$ws.Range("R3").Value = 2677.320694654988
$ws.Range("S3").Value = 0.1093313184008472
$ws.Range("T3").Value = 0.1093313184008472
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.824961333333333
$ws.Range("H4").Value = 5.474884
$ws.Range("I4").Value = 0.377840167393297
$ws.Range("J4").Value = 0.3778401673932969
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 302.9317738335067
$ws.Range("R4").Value = 2726.38596450156
$ws.Range("S4").Value = 0.1113349523512845
$ws.Range("T4").Value = 0.1113349523512845
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.824961333333333
$ws.Range("H5").Value = 5.474884
$ws.Range("I5").Value = 0.377840167393297
$ws.Range("J5").Value = 0.3778401673932969
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 120.8612293912725
$ws.Range("R5").Value = 1087.751064521452
$ws.Range("S5").Value = 0.04441950425045388
$ws.Range("T5").Value = 0.04441950425045387
$ws.Range("I6").Value = 0.03077064395059555
$ws.Range("J6").Value = 0.03077064395059554
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 24.98476272132556
$ws.Range("R6").Value = 224.86286449193
$ws.Range("S6").Value = 0.00918252097456029
$ws.Range("T6").Value = 0.00918252097456029
$ws.Range("I7").Value = 0.03077064395059555
$ws.Range("J7").Value = 0.03077064395059554
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("S7").Value = 0.008903751801644334
$ws.Range("T7").Value = 0.008903751801644332
$ws.Range("I8").Value = 0.03077064395059555
$ws.Range("J8").Value = 0.03077064395059554
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 24.67023508448333
$ws.Range("R8").Value = 222.03211576035
$ws.Range("S8").Value = 0.009066924254487484
$ws.Range("T8").Value = 0.009066924254487482
$ws.Range("I9").Value = 0.03077064395059555
$ws.Range("J9").Value = 0.03077064395059554
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 9.842727634510558
$ws.Range("R9").Value = 88.584548710595
$ws.Range("S9").Value = 0.003617446919903439
$ws.Range("T9").Value = 0.003617446919903438
$ws.Range("G10").Value = 2.658767
$ws.Range("H10").Value = 7.976300999999999
$ws.Range("I10").Value = 0.5504713716344166
$ws.Range("J10").Value = 0.5504713716344165
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 446.9648612895646
$ws.Range("R10").Value = 4022.683751606081
$ws.Range("S10").Value = 0.1642706900786252
$ws.Range("T10").Value = 0.1642706900786252
$ws.Range("G11").Value = 2.658767
$ws.Range("H11").Value = 7.976300999999999
$ws.Range("I11").Value = 0.5504713716344166
$ws.Range("J11").Value = 0.5504713716344165
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("Q11").Value = 433.3956001847563
$ws.Range("R11").Value = 3900.560401662807
$ws.Range("S11").Value = 0.1592836495333957
$ws.Range("T11").Value = 0.1592836495333957
$ws.Range("G12").Value = 2.658767
$ws.Range("H12").Value = 7.976300999999999
$ws.Range("I12").Value = 0.5504713716344166
$ws.Range("J12").Value = 0.5504713716344165
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 441.3381197775099
$ws.Range("R12").Value = 3972.04307799759
$ws.Range("S12").Value = 0.1622027227927574
$ws.Range("T12").Value = 0.1622027227927573
$ws.Range("G13").Value = 2.658767
$ws.Range("H13").Value = 7.976300999999999
$ws.Range("I13").Value = 0.5504713716344166
$ws.Range("J13").Value = 0.5504713716344165
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 176.0814557632337
$ws.Range("R13").Value = 1584.733101869103
$ws.Range("S13").Value = 0.06471430922963838
$ws.Range("T13").Value = 0.06471430922963836
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.1976323333333333
$ws.Range("H14").Value = 0.592897
$ws.Range("I14").Value = 0.04091781702169097
$ws.Range("J14").Value = 0.04091781702169097
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 33.22393743215044
$ws.Range("R14").Value = 299.015436889354
$ws.Range("S14").Value = 0.01221062235935512
$ws.Range("T14").Value = 0.01221062235935512
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.1976323333333333
$ws.Range("H15").Value = 0.592897
$ws.Range("I15").Value = 0.04091781702169097
$ws.Range("J15").Value = 0.04091781702169097
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 32.21530270268656
$ws.Range("R15").Value = 289.937724324179
$ws.Range("S15").Value = 0.01183992403965218
$ws.Range("T15").Value = 0.01183992403965218
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.1976323333333333
$ws.Range("H16").Value = 0.592897
$ws.Range("I16").Value = 0.04091781702169097
$ws.Range("J16").Value = 0.04091781702169097
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 32.80568865213667
$ws.Range("R16").Value = 295.25119786923
$ws.Range("S16").Value = 0.01205690554251369
$ws.Range("T16").Value = 0.01205690554251369
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.1976323333333333
$ws.Range("H17").Value = 0.592897
$ws.Range("I17").Value = 0.04091781702169097
$ws.Range("J17").Value = 0.04091781702169097
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 13.08854403534345
$ws.Range("R17").Value = 117.796896318091
$ws.Range("S17").Value = 0.004810365080169981
$ws.Range("T17").Value = 0.004810365080169981
